{"js": "// Fill in the header row of the report table and update the two data\n// rows to reflect the new inventory entries (Report Header table in\n// Django action).\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount,values\");\nawait context.sync();\n\n// New contents for the first three rows of the (single) table:\n//  - row 0: column headers (previously empty cells)\n//  - row 1: id=3, \u041f\u042d\u0412\u041c, inv=122, \u0432/\u0447 1234\n//  - row 2: id=1, \u0421\u0435\u0440\u0432\u0435\u0440, inv=\u0418\u0422111123, \u0413\u041f\u041a\nconst newValues = [\n  [\"ID\", \"\u041d\u0430\u0438\u043c\u0435\u043d\u043e\u0432\u0430\u043d\u0438\u0435 \u0442\u0435\u0445\u043d\u0438\u043a\u0438\", \"\u0418\u043d\u0432\u0435\u043d\u0442\u0430\u0440\u043d\u044b\u0439 \u043d\u043e\u043c\u0435\u0440\", \"\u041f\u043e\u0434\u0440\u0430\u0437\u0434\u0435\u043b\u0435\u043d\u0438\u0435\"],\n  [\"3\", \"\u041f\u042d\u0412\u041c\", \"122\", \"\u0432/\u0447 1234\"],\n  [\"1\", \"\u0421\u0435\u0440\u0432\u0435\u0440\", \"\u0418\u0422111123\", \"\u0413\u041f\u041a\"],\n];\n\nfor (let r = 0; r < newValues.length && r < table.rowCount; r++) {\n  for (let c = 0; c < newValues[r].length; c++) {\n    const cell = table.getCell(r, c);\n    cell.body.insertText(newValues[r][c], Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Fill in the header row of the report table and update the two data\n# rows to reflect the new inventory entries (Report Header table in\n# Django action).\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Row 1: previously-empty header cells -> column titles.\n$t.Cell(1, 1).Range.Text = \"ID\"\n$t.Cell(1, 2).Range.Text = \"\u041d\u0430\u0438\u043c\u0435\u043d\u043e\u0432\u0430\u043d\u0438\u0435 \u0442\u0435\u0445\u043d\u0438\u043a\u0438\"\n$t.Cell(1, 3).Range.Text = \"\u0418\u043d\u0432\u0435\u043d\u0442\u0430\u0440\u043d\u044b\u0439 \u043d\u043e\u043c\u0435\u0440\"\n$t.Cell(1, 4).Range.Text = \"\u041f\u043e\u0434\u0440\u0430\u0437\u0434\u0435\u043b\u0435\u043d\u0438\u0435\"\n\n# Row 2: id 4 -> 3, inventory 555666 -> 122, unit \u0432/\u0447 2044 -> \u0432/\u0447 1234.\n$t.Cell(2, 1).Range.Text = \"3\"\n$t.Cell(2, 3).Range.Text = \"122\"\n$t.Cell(2, 4).Range.Text = \"\u0432/\u0447 1234\"\n\n# Row 3: id 2 -> 1, name \u041c\u0424\u0423 \u04103 -> \u0421\u0435\u0440\u0432\u0435\u0440, inventory 456775 -> \u0418\u0422111123,\n# unit \u0432/\u0447 2044 -> \u0413\u041f\u041a.\n$t.Cell(3, 1).Range.Text = \"1\"\n$t.Cell(3, 2).Range.Text = \"\u0421\u0435\u0440\u0432\u0435\u0440\"\n$t.Cell(3, 3).Range.Text = \"\u0418\u0422111123\"\n$t.Cell(3, 4).Range.Text = \"\u0413\u041f\u041a\"\n"}
